$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update availability values: mark additional slots as fully available (3)
$ws.Range("C4").Value = 3
$ws.Range("P4").Value = 3
$ws.Range("C5").Value = 3
$ws.Range("P5").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("P6").Value = 3
$ws.Range("C7").Value = 3
